$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.439.90"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.916.42"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.70%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.33"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4822"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4077"
$ws.Range("E8").Value = "  +0.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08235"
$ws.Range("E9").Value = "  +2.75%  "
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.48"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.921.56"
$ws.Range("E12").Value = "  +5.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.072"
$ws.Range("E13").Value = "  +2.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.242"
$ws.Range("E14").Value = "  +2.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.30"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06808"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.68"
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.007"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.458.81"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.656"
$ws.Range("E22").Value = "  +2.54%  "
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.183.77"
$ws.Range("E25").Value = "  +5.08%  "
$ws.Range("E26").Value = "  +9.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.86"
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.114"
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.39"
$ws.Range("E30").Value = "  +2.23%  "
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09575"
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.682"
$ws.Range("E33").Value = "  +6.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.550"
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.375"
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02290"
$ws.Range("E36").Value = "  +1.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06109"
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.180"
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5990"
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.054"
$ws.Range("E40").Value = "  +3.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.81"
$ws.Range("E41").Value = "  +7.10%  "
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.428"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07631"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.40"
$ws.Range("E46").Value = "  +1.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5590"
$ws.Range("E47").Value = "  +1.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.956"
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.72"
$ws.Range("E49").Value = "  +4.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.426"
$ws.Range("E50").Value = "  +4.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.43"
$ws.Range("E51").Value = "  +1.17%  "
